# Update recomputed TPM-derived NATMI edge statistics for the Bmp10-Acvrl1 pair.
# Column A (sending cluster) and E:T (expression/specificity stats) are refreshed
# per updated TPM input; column D (target cluster) values are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row=2; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=41.474781; N=82.949562; O=0.6394851352970483; P=0.5499336598827257; Q=15.240544225092; R=91.443265350552; S=0.3563603059139592; T=0.3064567359758236 },
    @{ Row=3; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=9.567994999999998; N=28.703985; O=0.1475255668522151; P=0.1902998297238611; Q=3.515906472006666; R=31.64315824805999; S=0.08221028641919126; T=0.1060467269567861 },
    @{ Row=4; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=5.805205666666667; N=17.415617; O=0.0895084348046473; P=0.1154609351153152; Q=2.133211835370222; R=19.198906518332; S=0.04987958507283698; T=0.06434190865076618 },
    @{ Row=5; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=2.2592025; N=4.518405; O=0.03483385280264482; P=0.02995583024877705; Q=0.8301785997300001; R=4.98107159838; S=0.01941155744792436; T=0.01669322434899465 },
    @{ Row=6; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=1.733200666666667; N=5.199602000000001; O=0.02672361459413777; P=0.03447198621487044; Q=0.6368911607102222; R=5.732020446392; S=0.01489203571161983; T=0.01920990320953551 },
    @{ Row=7; A="Neutrophils"; E=1; F=0.3333333333333333; G=0.3674653333333333; H=1.102396; I=0.5572612813719676; J=0.5572612813719677; M=4.016136; N=12.048408; O=0.06192339564930666; P=0.07987775881445054; Q=1.475790753952; R=13.282116785568; S=0.03450751080643596; T=0.0445127822300617 },
    @{ Row=8; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=41.474781; N=82.949562; O=0.6394851352970483; P=0.5499336598827257; Q=12.108465538461; R=72.650793230766; S=0.283124829383089; T=0.2434769239069021 },
    @{ Row=9; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=9.567994999999998; N=28.703985; O=0.1475255668522151; P=0.1902998297238611; Q=2.793353814928333; R=25.140184334355; S=0.06531528043302383; T=0.08425310276707501 },
    @{ Row=10; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=5.805205666666667; N=17.415617; O=0.0895084348046473; P=0.1154609351153152; Q=1.694816248903444; R=15.253346240131; S=0.03962884973181031; T=0.05111902646454905 },
    @{ Row=11; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=2.2592025; N=4.518405; O=0.03483385280264482; P=0.02995583024877705; Q=0.6595688984025001; R=3.957413390415001; S=0.01542229535472046; T=0.0132626058997824 },
    @{ Row=12; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=1.733200666666667; N=5.199602000000001; O=0.02672361459413777; P=0.03447198621487044; Q=0.5060038904984445; R=4.554035014486001; S=0.01183157888251794; T=0.01526208300533493 },
    @{ Row=13; A="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.2919476666666667; H=0.875843; I=0.4427387186280323; J=0.4427387186280324; M=4.016136; N=12.048408; O=0.06192339564930666; P=0.07987775881445054; Q=1.172501534216; R=10.552513807944; S=0.0274158848428707; T=0.03536497658438885 }
)

foreach ($r in $rowData) {
    foreach ($col in @("A","E","F","G","H","I","J","M","N","O","P","Q","R","S","T")) {
        $ws.Range("$col$($r.Row)").Value = $r[$col]
    }
}
